$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26. This shifts the blank gap row and the
# "Total" row down by one (old row 27 -> row 28) and automatically adjusts
# the relative reference in the SUM formula from C2:C26 to C2:C27.
$ws.Rows.Item(26).Insert()

# Fill in the new timesheet entry in the freshly inserted row 26.
$ws.Range("A26").Value = "Made corrections and added one graphic"
$ws.Range("B26").Value = 43346
$ws.Range("C26").Value = 0.5

# The sheet had scrolled down (topLeftCell) and had C26 selected; after the
# edit the view resets to the top and the new Total cell (C27) is selected.
$ws.Application.ActiveWindow.ScrollRow = 1
$null = $ws.Range("C27").Select()
